# Udated cure data (v44)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 73
$colCount = 2
$values = New-Object 'object[,]' $rowCount, $colCount

$values[0,0] = 'Ziekenhuisorganisatie'; $values[0,1] = 'Status inleveren routekaart'
$values[1,0] = 'Admiraal De Ruyter Ziekenhuis'; $values[1,1] = 'definitief en/of vastgesteld RvB'
$values[2,0] = 'Albert Schweitzer ziekenhuis'; $values[2,1] = 'definitief en/of vastgesteld RvB'
$values[3,0] = 'Alrijne Zorggroep'; $values[3,1] = 'voorlopig'
$values[4,0] = 'Amphia Ziekenhuis'; $values[4,1] = 'definitief en/of vastgesteld RvB'
$values[5,0] = 'Amsterdam UMC'; $values[5,1] = 'definitief en/of vastgesteld RvB'
$values[6,0] = 'Antoni van Leeuwenhoek'; $values[6,1] = 'definitief en/of vastgesteld RvB'
$values[7,0] = 'Antonius Zorggroep'; $values[7,1] = 'definitief en/of vastgesteld RvB'
$values[8,0] = 'Bravis ziekenhuis'; $values[8,1] = 'definitief en/of vastgesteld RvB'
$values[9,0] = 'Canisius-Wilhelmina Ziekenhuis'; $values[9,1] = 'definitief en/of vastgesteld RvB'
$values[10,0] = 'Catharina Ziekenhuis'; $values[10,1] = 'definitief en/of vastgesteld RvB'
$values[11,0] = 'CuraMare'; $values[11,1] = 'voorlopig'
$values[12,0] = 'Deventer Ziekenhuis'; $values[12,1] = 'definitief en/of vastgesteld RvB'
$values[13,0] = 'Diakonessenhuis'; $values[13,1] = 'definitief en/of vastgesteld RvB'
$values[14,0] = 'Elisabeth-TweeSteden Ziekenhuis'; $values[14,1] = 'definitief en/of vastgesteld RvB'
$values[15,0] = 'Elkerliek Ziekenhuis'; $values[15,1] = 'definitief en/of vastgesteld RvB'
$values[16,0] = 'Erasmus MC'; $values[16,1] = 'definitief en/of vastgesteld RvB'
$values[17,0] = 'Flevoziekenhuis'; $values[17,1] = 'definitief en/of vastgesteld RvB'
$values[18,0] = 'Franciscus Gasthuis & Vlietland Groep'; $values[18,1] = 'definitief en/of vastgesteld RvB'
$values[19,0] = 'Gelre Ziekenhuizen'; $values[19,1] = 'definitief en/of vastgesteld RvB'
$values[20,0] = 'Groene Hart Ziekenhuis'; $values[20,1] = 'definitief en/of vastgesteld RvB'
$values[21,0] = 'Haaglanden Medisch Centrum'; $values[21,1] = 'voorlopig'
$values[22,0] = 'IJsselland Ziekenhuis'; $values[22,1] = 'definitief en/of vastgesteld RvB'
$values[23,0] = 'Isala Klinieken'; $values[23,1] = 'definitief en/of vastgesteld RvB'
$values[24,0] = 'Jeroen Bosch Ziekenhuis'; $values[24,1] = 'definitief en/of vastgesteld RvB'
$values[25,0] = 'LUMC'; $values[25,1] = 'definitief en/of vastgesteld RvB'
$values[26,0] = 'LangeLand Ziekenhuis'; $values[26,1] = 'definitief en/of vastgesteld RvB'
$values[27,0] = 'Laurentius Ziekenhuis'; $values[27,1] = 'definitief en/of vastgesteld RvB'
$values[28,0] = 'Maasstad ziekenhuis'; $values[28,1] = 'definitief en/of vastgesteld RvB'
$values[29,0] = 'Maastricht UMC+'; $values[29,1] = 'definitief en/of vastgesteld RvB'
$values[30,0] = 'Martini Ziekenhuis'; $values[30,1] = 'definitief en/of vastgesteld RvB'
$values[31,0] = 'Meander Medisch Centrum'; $values[31,1] = 'definitief en/of vastgesteld RvB'
$values[32,0] = 'Medisch Centrum Leeuwarden'; $values[32,1] = 'definitief en/of vastgesteld RvB'
$values[33,0] = 'Medisch Spectrum Twente'; $values[33,1] = 'definitief en/of vastgesteld RvB'
$values[34,0] = 'Máxima Medisch Centrum'; $values[34,1] = 'definitief en/of vastgesteld RvB'
$values[35,0] = 'Nij Smellinghe'; $values[35,1] = 'definitief en/of vastgesteld RvB'
$values[36,0] = 'Noordwest Ziekenhuisgroep'; $values[36,1] = 'definitief en/of vastgesteld RvB'
$values[37,0] = 'Ommelander Ziekenhuis Groep'; $values[37,1] = 'definitief en/of vastgesteld RvB'
$values[38,0] = 'Onze Lieve Vrouwe Gasthuis'; $values[38,1] = 'definitief en/of vastgesteld RvB'
$values[39,0] = 'Pantein'; $values[39,1] = 'definitief en/of vastgesteld RvB'
$values[40,0] = 'Prinses Máxima Centrum'; $values[40,1] = 'definitief en/of vastgesteld RvB'
$values[41,0] = 'Protestants Christelijk Ziekenhuis Ikazia'; $values[41,1] = 'definitief en/of vastgesteld RvB'
$values[42,0] = 'Radboudumc'; $values[42,1] = 'definitief en/of vastgesteld RvB'
$values[43,0] = 'Radiotherapiegroep'; $values[43,1] = 'voorlopig'
$values[44,0] = 'Rivas Zorggroep'; $values[44,1] = 'definitief en/of vastgesteld RvB'
$values[45,0] = 'Rode Kruis Ziekenhuis'; $values[45,1] = 'voorlopig'
$values[46,0] = 'Slingeland Ziekenhuis (Santiz)'; $values[46,1] = 'definitief en/of vastgesteld RvB'
$values[47,0] = 'Spaarne Gasthuis'; $values[47,1] = 'definitief en/of vastgesteld RvB'
$values[48,0] = 'Spijkenisse Medisch Centrum'; $values[48,1] = 'definitief en/of vastgesteld RvB'
$values[49,0] = 'St Jansdal'; $values[49,1] = 'definitief en/of vastgesteld RvB'
$values[50,0] = 'St. Antonius Ziekenhuis'; $values[50,1] = 'definitief en/of vastgesteld RvB'
$values[51,0] = 'St. Jans Gasthuis'; $values[51,1] = 'definitief en/of vastgesteld RvB'
$values[52,0] = 'Stichting BovenIJ ziekenhuis'; $values[52,1] = 'definitief en/of vastgesteld RvB'
$values[53,0] = 'Stichting Dijklander Ziekenhuis'; $values[53,1] = 'definitief en/of vastgesteld RvB'
$values[54,0] = 'Stichting Reinier Haga Groep'; $values[54,1] = 'definitief en/of vastgesteld RvB'
$values[55,0] = 'Stichting Reinier de Graaf Groep'; $values[55,1] = 'definitief en/of vastgesteld RvB'
$values[56,0] = 'Stichting Rijnstate Ziekenhuis'; $values[56,1] = 'definitief en/of vastgesteld RvB'
$values[57,0] = 'Stichting St. Anna Zorggroep'; $values[57,1] = 'definitief en/of vastgesteld RvB'
$values[58,0] = 'Stichting Tergooi'; $values[58,1] = 'definitief en/of vastgesteld RvB'
$values[59,0] = 'Stichting VieCuri-Vitaal'; $values[59,1] = 'definitief en/of vastgesteld RvB'
$values[60,0] = 'Stichting Ziekenhuisgroep Twente (ZGT)'; $values[60,1] = 'definitief en/of vastgesteld RvB'
$values[61,0] = 'Streekziekenhuis Koningin Beatrix (Santiz)'; $values[61,1] = 'definitief en/of vastgesteld RvB'
$values[62,0] = 'Tjongerschans'; $values[62,1] = 'definitief en/of vastgesteld RvB'
$values[63,0] = 'Treant Zorggroep'; $values[63,1] = 'definitief en/of vastgesteld RvB'
$values[64,0] = 'UMC Utrecht'; $values[64,1] = 'definitief en/of vastgesteld RvB'
$values[65,0] = 'UMCG'; $values[65,1] = 'definitief en/of vastgesteld RvB'
$values[66,0] = 'Wilhelmina Ziekenhuis Assen'; $values[66,1] = 'definitief en/of vastgesteld RvB'
$values[67,0] = 'Ziekenhuis Amstelland'; $values[67,1] = 'voorlopig'
$values[68,0] = 'Ziekenhuis Bernhoven'; $values[68,1] = 'definitief en/of vastgesteld RvB'
$values[69,0] = 'Ziekenhuis De Gelderse Vallei'; $values[69,1] = 'definitief en/of vastgesteld RvB'
$values[70,0] = 'Ziekenhuis Rivierenland'; $values[70,1] = 'definitief en/of vastgesteld RvB'
$values[71,0] = 'ZorgSaam Zeeuws-Vlaanderen'; $values[71,1] = 'definitief en/of vastgesteld RvB'
$values[72,0] = 'Zuyderland Ziekenhuis'; $values[72,1] = 'definitief en/of vastgesteld RvB'

$target = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($rowCount, $colCount))
$target.Value = $values
